# Apply scraped-symbol-list update (Mon Jan 23 22:34:50 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'306.05"
$ws.Range('E2').Value = "'1.63%"
$ws.Range('E3').Value = "'-0.64%"
$ws.Range('D4').Value = "'5.067"
$ws.Range('E4').Value = "'1.42%"
$ws.Range('D5').Value = "'0.07921"
$ws.Range('D6').Value = "'2.203"
$ws.Range('E6').Value = "'5.77%"
$ws.Range('B7').Value = "'KuCoinToken"
$ws.Range('C7').Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range('D7').Value = "'8.017"
$ws.Range('E7').Value = "'1.30%"
$ws.Range('B8').Value = "'MXToken"
$ws.Range('C8').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D8').Value = "'0.9307"
$ws.Range('E8').Value = "'1.70%"
$ws.Range('B9').Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range('C9').Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range('D9').Value = "'0.09843"
$ws.Range('E9').Value = "'1.79%"
$ws.Range('B10').Value = "'WazirX"
$ws.Range('C10').Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range('D10').Value = "'0.1877"
$ws.Range('E10').Value = "'0.67%"
$ws.Range('B11').Value = "'MandalaExchangeToken"
$ws.Range('C11').Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range('D11').Value = "'0.09113"
$ws.Range('E11').Value = "'6.90%"
$ws.Range('B12').Value = "'BitrueCoin"
$ws.Range('C12').Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range('D12').Value = "'0.03685"
$ws.Range('E12').Value = "'4.14%"
$ws.Range('B13').Value = "'BitMartToken"
$ws.Range('C13').Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range('D13').Value = "'0.09917"
$ws.Range('E13').Value = "'-0.44%"
$ws.Range('B14').Value = "'BitForexToken"
$ws.Range('C14').Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range('D14').Value = "'0.001435"
$ws.Range('E14').Value = "'-3.18%"
$ws.Range('B15').Value = "'TigerCash"
$ws.Range('C15').Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range('D15').Value = "'0.005621"
$ws.Range('E15').Value = "'-0.96%"
$ws.Range('B16').Value = "'LEO"
$ws.Range('C16').Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range('D16').Value = "'3.465"
$ws.Range('E16').Value = "'-0.01%"
$ws.Range('B17').Value = "'GateToken"
$ws.Range('C17').Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range('D17').Value = "'4.176"
$ws.Range('E17').Value = "'3.63%"
$ws.Range('D20').Value = "'0.1346"
$ws.Range('E20').Value = "'1.29%"
$ws.Range('D21').Value = "'5.101"
$ws.Range('E21').Value = "'7.24%"
$ws.Range('D22').Value = "'0.2190"
$ws.Range('E22').Value = "'-0.49%"
$ws.Range('D23').Value = "'0.04559"
$ws.Range('E23').Value = "'-0.55%"
$ws.Range('D24').Value = "'0.001240"
$ws.Range('E24').Value = "'0.56%"
$ws.Range('D25').Value = "'0.004786"
$ws.Range('E25').Value = "'-6.03%"
$ws.Range('D26').Value = "'0.0001299"
$ws.Range('E26').Value = "'-7.30%"
$ws.Range('D39').Value = "'0.01931"
$ws.Range('E39').Value = "'9.96%"
$ws.Range('E40').Value = "'6.84%"
$ws.Range('D41').Value = "'0.007811"
$ws.Range('E41').Value = "'4.33%"
$ws.Range('D42').Value = "'0.1395"
$ws.Range('E42').Value = "'0.31%"
$ws.Range('D43').Value = "'0.007804"
$ws.Range('E43').Value = "'1.00%"
$ws.Range('D44').Value = "'0.002110"
$ws.Range('E44').Value = "'-5.89%"
$ws.Range('D45').Value = "'0.01123"
$ws.Range('E45').Value = "'9.22%"
$ws.Range('D46').Value = "'0.00006242"
$ws.Range('E46').Value = "'-0.77%"
$ws.Range('E47').Value = "'-0.17%"
$ws.Range('D48').Value = "'51.90"
$ws.Range('E48').Value = "'36.40%"
$ws.Range('E49').Value = "'-10.15%"
$ws.Range('E50').Value = "'-0.17%"
$ws.Range('E51').Value = "'-0.17%"
